$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Labels in column C
$ws.Range("C1").Value = "average"
$ws.Range("C2").Value = "max"
$ws.Range("C3").Value = "min"

# Formulas in column D
$ws.Range("D1").Formula = "=AVERAGE(A:A)"
$ws.Range("D2").Formula = "=MAX(A:A)"
$ws.Range("D3").Formula = "=MIN(A:A)"

# Bold the average result cell
$ws.Range("D1").Font.Bold = $true

# Set column D width to best-fit the longest entry (mirrors Excel's
# auto-fit-on-entry behavior for the new MAX/MIN/AVERAGE column)
$ws.Columns.Item(4).ColumnWidth = 10.75

# Move selection
$ws.Range("G20").Select()
